$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text (e.g. "1.001") must be forced to
# text first, otherwise assigning .Value lets Excel coerce them to real
# numbers (dropping significant trailing zeros, e.g. "6.340" -> 6.34).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.596.56"
$ws.Range("D3").Value = "1.829.18"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "317.75"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "0.5406"
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").Value = "0.4018"
$ws.Range("E8").Value = "  +6.56%  "
$ws.Range("D9").Value = "0.07744"
$ws.Range("E9").Value = "  +4.37%  "
$ws.Range("D10").Value = "1.121"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "41.97"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "21.45"
$ws.Range("E12").Value = "  +4.65%  "
$ws.Range("D13").Value = "6.340"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "7.660"
$ws.Range("E14").Value = "  +6.01%  "
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "1.830.27"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("E18").Value = "  +1.35%  "
$ws.Range("D19").Value = "0.06595"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "17.76"
$ws.Range("E20").Value = "  +3.21%  "
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "6.083"
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("D23").Value = "28.590.54"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").Value = "11.22"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  +8.23%  "
$ws.Range("D26").Value = "2.471"
$ws.Range("E26").Value = "  +8.23%  "
$ws.Range("D27").Value = "158.31"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").Value = "20.83"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("D29").Value = "2.038.95"
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("D30").Value = "124.43"
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").Value = "0.1116"
$ws.Range("E32").Value = "  +5.30%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.692"
$ws.Range("E33").Value = "  +2.82%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.07496"
$ws.Range("E34").Value = "  +16.31%  "
$ws.Range("D35").Value = "3.650"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").Value = "0.2259"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").Value = "0.02362"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").Value = "8.950"
$ws.Range("E38").Value = "  +5.91%  "
$ws.Range("D39").Value = "5.227"
$ws.Range("E39").Value = "  +4.72%  "
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "11.41"
$ws.Range("E40").Value = "  +2.93%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6319"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").Value = "1.192"
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "1.405"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("D46").Value = "0.5902"
$ws.Range("E46").Value = "  +2.31%  "
$ws.Range("D47").Value = "3.712"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").Value = "125.58"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "2.004"
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "0.06914"
$ws.Range("E51").Value = "  +1.59%  "

# Restore default styling on those cells (source data carries no explicit
# cell style), now that the text values are safely stored.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
